# Update odds values in the active worksheet to reflect the latest
# FlashScore data pull, per the commit "Atualizando o arquivo XLSX".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 ---
$ws.Range("Q4").Value = 2.6
$ws.Range("R4").Value = 1.48

# --- Row 8 ---
$ws.Range("G8").Value = 3.1
$ws.Range("H8").Value = 3.1

# --- Row 18 ---
$ws.Range("G18").Value = 1.24
$ws.Range("H18").Value = 5.6
$ws.Range("I18").Value = 10.75
$ws.Range("J18").Value = 1.62
$ws.Range("K18").Value = 2.7
$ws.Range("L18").Value = 8.25
$ws.Range("N18").Value = 9.5
$ws.Range("O18").Value = 1.16
$ws.Range("P18").Value = 4.65
$ws.Range("Q18").Value = 1.47
$ws.Range("R18").Value = 2.5
$ws.Range("S18").Value = 1.26
$ws.Range("T18").Value = 3.5
$ws.Range("U18").Value = 1.98
$ws.Range("V18").Value = 1.75
$ws.Range("W18").Value = 8.25
$ws.Range("Y18").Value = 9
$ws.Range("Z18").Value = 7.5
$ws.Range("AB18").Value = 27
$ws.Range("AC18").Value = 9.5
$ws.Range("AD18").Value = 11.5
$ws.Range("AE18").Value = 23
$ws.Range("AH18").Value = 30
$ws.Range("AI18").Value = 90
$ws.Range("AJ18").Value = 32
$ws.Range("AK18").Value = 350
$ws.Range("AL18").Value = 120
$ws.Range("AM18").Value = 100
$ws.Range("AO18").Value = 5.2
$ws.Range("AP18").Value = 14.5
$ws.Range("AQ18").Value = 12.5
$ws.Range("AR18").Value = 35
$ws.Range("AT18").Value = 3.5
$ws.Range("AX18").Value = 10.75
$ws.Range("AY18").Value = 60
$ws.Range("AZ18").Value = 50
$ws.Range("BA18").Value = 450
$ws.Range("BB18").Value = 400
